$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear out the old 6-column / 5-row table so we can rebuild it ---
$ws.Range("A1:F5").Clear()

# --- Header row ---
$ws.Range("A1").Value = "Combination"
$ws.Range("B1").Value = "Number_Of_Records"
$ws.Range("C1").Value = "Candidate_Pairs"
$ws.Range("D1").Value = "Time_Required (Mins)"
$ws.Range("E1").Value = "Time_Required (Sec)"

# --- Row 2 ---
$ws.Range("A2").Value = "Single Country"
$ws.Range("B2").Value = 3500
$ws.Range("C2").Formula = "=B2*(B2-1)/2"
$ws.Range("D2").Value = 2.2000000000000002
$ws.Range("E2").Formula = "=D2*60"

# --- Row 3 ---
$ws.Range("A3").Value = "Single Country"
$ws.Range("B3").Value = 5000
$ws.Range("C3").Formula = "=B3*(B3-1)/2"
$ws.Range("D3").Formula = "=(C3/C2)*D2"
$ws.Range("E3").Formula = "=D3*60"

# --- Row 4 ---
$ws.Range("A4").Value = "Single Country"
$ws.Range("B4").Value = 20398
$ws.Range("C4").Formula = "=B4*(B4-1)/2"
$ws.Range("D4").Value = "RAM usage crashes code"
$ws.Range("E4").Value = "RAM usage crashes code"

# --- Formatting: thin box border around the whole table, and a thousands
#     number format on B4 (20,398) ---
$ws.Range("A1:E4").Borders.LineStyle = 1
$ws.Range("B4").NumberFormat = "#,##0"

# --- Column widths for the now-5-wide table ---
$ws.Columns.Item(4).ColumnWidth = 18.619791666666668
$ws.Columns.Item(5).ColumnWidth = 20.709635416666668

# --- Selection, matching the saved workbook state ---
[void]$ws.Range("A6").Select()
